$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Heading: "The holographic principle and gravitational degeneracy"
# becomes 4 runs: "The holographic principle and " / "the " /
# "gravitational " / "keystone"
$headingXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="31899214" w14:textId="104FD7E0" w:rsidR="00900FED" w:rsidRDefault="00000000"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The holographic principle and </w:t></w:r><w:r><w:t xml:space="preserve">the </w:t></w:r><w:r><w:t xml:space="preserve">gravitational </w:t></w:r><w:r><w:t>keystone</w:t></w:r></w:p>
'@

$headingRng = $d.Content
$found = $headingRng.Find.Execute("The holographic principle and gravitational degeneracy")
if ($found) {
    $para = $headingRng.Paragraphs(1)
    $paraRng = $para.Range
    $paraRng.Collapse(0)
    $paraRng.InsertXML($headingXml)
}

# --- Change 2 -------------------------------------------------------------
# "It's a matter of *degeneracy*, and minimum size ..." loses the italic
# "degeneracy" aside, and later "... each oscillator acts as a keystone,
# stopping one another from falling further toward the ..." gains an
# inline <n> oscillator count (OMML) and gets re-wrapped into fewer runs.
$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0293BE3C" w14:textId="67A41C47" w:rsidR="00900FED" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">It’s a matter of minimum size – there is no </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>singularity</w:t></w:r><w:r><w:t xml:space="preserve"> in this model of the black hole process.</w:t></w:r><w:r w:rsidR="00601E9A"><w:t xml:space="preserve"> In effect, each</w:t></w:r><w:r><w:t xml:space="preserve"> of the </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>n</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> oscillator acts as a keystone</w:t></w:r><w:r w:rsidR="00143CEB"><w:t xml:space="preserve">, stopping one another from falling further toward the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00143CEB"><w:t>centre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00143CEB"><w:t xml:space="preserve"> of the black hole process.</w:t></w:r></w:p>
'@

$bodyRng = $d.Content
$found2 = $bodyRng.Find.Execute("It’s a matter of")
if (-not $found2) {
    $bodyRng = $d.Content
    $found2 = $bodyRng.Find.Execute("It's a matter of")
}
if ($found2) {
    $bodyPara = $bodyRng.Paragraphs(1)
    $bodyParaRng = $bodyPara.Range
    $bodyParaRng.Collapse(0)
    $bodyParaRng.InsertXML($bodyXml)
}
